$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value2 = "Datos actualizados a 14 de Abril de 2020 a las 13:22"

# Rows 11-11
$arr11 = New-Object 'object[,]' 1,8
$arr11[0,0] = "Iran"
$arr11[0,1] = 74877
$arr11[0,2] = 1574
$arr11[0,3] = 48129
$arr11[0,4] = 22065
$arr11[0,5] = 3691
$arr11[0,6] = 98
$arr11[0,7] = 4683
$ws.Range("A11:H11").Value2 = $arr11

# Rows 15-15
$arr15 = New-Object 'object[,]' 1,8
$arr15[0,0] = "Suiza"
$arr15[0,1] = 25807
$arr15[0,2] = 119
$arr15[0,3] = 13700
$arr15[0,4] = 10952
$arr15[0,5] = 386
$arr15[0,6] = 17
$arr15[0,7] = 1155
$ws.Range("A15:H15").Value2 = $arr15

# Rows 20-20
$arr20 = New-Object 'object[,]' 1,8
$arr20[0,0] = "Austria"
$arr20[0,1] = 14135
$arr20[0,2] = 94
$arr20[0,3] = 7633
$arr20[0,4] = 6118
$arr20[0,5] = 243
$arr20[0,6] = 16
$arr20[0,7] = 384
$ws.Range("A20:H20").Value2 = $arr20

# Rows 31-31
$arr31 = New-Object 'object[,]' 1,8
$arr31[0,0] = "Rumania"
$arr31[0,1] = 6879
$arr31[0,2] = 246
$arr31[0,3] = 1051
$arr31[0,4] = 5482
$arr31[0,5] = 241
$arr31[0,6] = 15
$arr31[0,7] = 346
$ws.Range("A31:H31").Value2 = $arr31

# Rows 63-67
$arr63 = New-Object 'object[,]' 5,8
$arr63[0,0] = "Barein"
$arr63[0,1] = 1522
$arr63[0,2] = 161
$arr63[0,3] = 645
$arr63[0,4] = 870
$arr63[0,5] = 3
$arr63[0,6] = 1
$arr63[0,7] = 7
$arr63[1,0] = "Hungria"
$arr63[1,1] = 1512
$arr63[1,2] = 54
$arr63[1,3] = 122
$arr63[1,4] = 1268
$arr63[1,5] = 58
$arr63[1,6] = 13
$arr63[1,7] = 122
$arr63[2,0] = "Irak"
$arr63[2,1] = 1378
$arr63[2,2] = 0
$arr63[2,3] = 717
$arr63[2,4] = 583
$arr63[2,5] = 0
$arr63[2,6] = 0
$arr63[2,7] = 78
$arr63[3,0] = "Estonia"
$arr63[3,1] = 1373
$arr63[3,2] = 41
$arr63[3,3] = 115
$arr63[3,4] = 1227
$arr63[3,5] = 11
$arr63[3,6] = 3
$arr63[3,7] = 31
$arr63[4,0] = "Nueva Zelanda"
$arr63[4,1] = 1366
$arr63[4,2] = 17
$arr63[4,3] = 628
$arr63[4,4] = 729
$arr63[4,5] = 4
$arr63[4,6] = 4
$arr63[4,7] = 9
$ws.Range("A63:H67").Value2 = $arr63

# Rows 70-70
$arr70 = New-Object 'object[,]' 1,8
$arr70[0,0] = "Kazajistan"
$arr70[0,1] = 1179
$arr70[0,2] = 88
$arr70[0,3] = 150
$arr70[0,4] = 1015
$arr70[0,5] = 21
$arr70[0,6] = 2
$arr70[0,7] = 14
$ws.Range("A70:H70").Value2 = $arr70

# Rows 72-75
$arr72 = New-Object 'object[,]' 4,8
$arr72[0,0] = "Uzbekistan"
$arr72[0,1] = 1113
$arr72[0,2] = 115
$arr72[0,3] = 99
$arr72[0,4] = 1010
$arr72[0,5] = 8
$arr72[0,6] = 0
$arr72[0,7] = 4
$arr72[1,0] = "Bosnia y Herzegovina"
$arr72[1,1] = 1080
$arr72[1,2] = 43
$arr72[1,3] = 218
$arr72[1,4] = 823
$arr72[1,5] = 4
$arr72[1,6] = 0
$arr72[1,7] = 39
$arr72[2,0] = "Lituania"
$arr72[2,1] = 1070
$arr72[2,2] = 8
$arr72[2,3] = 101
$arr72[2,4] = 945
$arr72[2,5] = 14
$arr72[2,6] = 0
$arr72[2,7] = 24
$arr72[3,0] = "Armenia"
$arr72[3,1] = 1067
$arr72[3,2] = 28
$arr72[3,3] = 265
$arr72[3,4] = 786
$arr72[3,5] = 30
$arr72[3,6] = 2
$arr72[3,7] = 16
$ws.Range("A72:H75").Value2 = $arr72

# Rows 100-103
$arr100 = New-Object 'object[,]' 4,8
$arr100[0,0] = "Malta"
$arr100[0,1] = 393
$arr100[0,2] = 9
$arr100[0,3] = 44
$arr100[0,4] = 346
$arr100[0,5] = 4
$arr100[0,6] = 0
$arr100[0,7] = 3
$arr100[1,0] = "Taiwan"
$arr100[1,1] = 393
$arr100[1,2] = 0
$arr100[1,3] = 124
$arr100[1,4] = 263
$arr100[1,5] = 0
$arr100[1,6] = 0
$arr100[1,7] = 6
$arr100[2,0] = "Reunion"
$arr100[2,1] = 391
$arr100[2,2] = 0
$arr100[2,3] = 40
$arr100[2,4] = 351
$arr100[2,5] = 3
$arr100[2,6] = 0
$arr100[2,7] = 0
$arr100[3,0] = "Jordania"
$arr100[3,1] = 391
$arr100[3,2] = 0
$arr100[3,3] = 215
$arr100[3,4] = 169
$arr100[3,5] = 5
$arr100[3,6] = 0
$arr100[3,7] = 7
$ws.Range("A100:H103").Value2 = $arr100

# Rows 114-114
$arr114 = New-Object 'object[,]' 1,8
$arr114[0,0] = "Vietnam"
$arr114[0,1] = 266
$arr114[0,2] = 1
$arr114[0,3] = 169
$arr114[0,4] = 97
$arr114[0,5] = 8
$arr114[0,6] = 0
$arr114[0,7] = 0
$ws.Range("A114:H114").Value2 = $arr114

# Rows 117-117
$arr117 = New-Object 'object[,]' 1,8
$arr117[0,0] = "Sri Lanka"
$arr117[0,1] = 219
$arr117[0,2] = 2
$arr117[0,3] = 59
$arr117[0,4] = 153
$arr117[0,5] = 1
$arr117[0,6] = 0
$arr117[0,7] = 7
$ws.Range("A117:H117").Value2 = $arr117

# Rows 133-133
$arr133 = New-Object 'object[,]' 1,8
$arr133[0,0] = "Madagascar"
$arr133[0,1] = 108
$arr133[0,2] = 2
$arr133[0,3] = 23
$arr133[0,4] = 85
$arr133[0,5] = 1
$arr133[0,6] = 0
$arr133[0,7] = 0
$ws.Range("A133:H133").Value2 = $arr133

# Rows 192-195
$arr192 = New-Object 'object[,]' 4,8
$arr192[0,0] = "Sierra Leona"
$arr192[0,1] = 11
$arr192[0,2] = 1
$arr192[0,3] = 0
$arr192[0,4] = 11
$arr192[0,5] = 0
$arr192[0,6] = 0
$arr192[0,7] = 0
$arr192[1,0] = "Montserrat"
$arr192[1,1] = 11
$arr192[1,2] = 0
$arr192[1,3] = 1
$arr192[1,4] = 10
$arr192[1,5] = 1
$arr192[1,6] = 0
$arr192[1,7] = 0
$arr192[2,0] = "Republica de Africa Central"
$arr192[2,1] = 11
$arr192[2,2] = 0
$arr192[2,3] = 3
$arr192[2,4] = 8
$arr192[2,5] = 0
$arr192[2,6] = 0
$arr192[2,7] = 0
$arr192[3,0] = "Groenlandia"
$arr192[3,1] = 11
$arr192[3,2] = 0
$arr192[3,3] = 11
$arr192[3,4] = 0
$arr192[3,5] = 0
$arr192[3,6] = 0
$arr192[3,7] = 0
$ws.Range("A192:H195").Value2 = $arr192

# Rows 206-207
$arr206 = New-Object 'object[,]' 2,8
$arr206[0,0] = "Islas Malvinas"
$arr206[0,1] = 5
$arr206[0,2] = 0
$arr206[0,3] = 1
$arr206[0,4] = 4
$arr206[0,5] = 0
$arr206[0,6] = 0
$arr206[0,7] = 0
$arr206[1,0] = "Burundi"
$arr206[1,1] = 5
$arr206[1,2] = 0
$arr206[1,3] = 0
$arr206[1,4] = 4
$arr206[1,5] = 0
$arr206[1,6] = 0
$arr206[1,7] = 1
$ws.Range("A206:H207").Value2 = $arr206
